$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price values in column D are stored as plain text in the source workbook
# (e.g. "1.005", "0.00001070"). Excel auto-parses single-dot numeric-looking
# strings as numbers on assignment, so those specific cells are pre-formatted
# as Text ("@") to preserve the literal text (matches/leading zeros/etc.).

# Row 2
$ws.Range("D2").Value = '27.467.98'
$ws.Range("E2").Value = '  -3.19%  '
# Row 3
$ws.Range("D3").Value = '1.754.66'
$ws.Range("E3").Value = '  -2.71%  '
# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.005'
$ws.Range("E4").Value = '  +0.26%  '
# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '324.27'
$ws.Range("E5").Value = '  -1.06%  '
# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.002'
$ws.Range("E6").Value = '  +0.24%  '
# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4416'
$ws.Range("E7").Value = '  -0.95%  '
# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3693'
$ws.Range("E8").Value = '  -2.21%  '
# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '44.87'
$ws.Range("E9").Value = '  +0.68%  '
# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.07638'
$ws.Range("E10").Value = '  +1.85%  '
# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.111'
$ws.Range("E11").Value = '  -3.25%  '
# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.003'
$ws.Range("E12").Value = '  +0.21%  '
# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '21.56'
$ws.Range("E13").Value = '  -4.50%  '
# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.148'
$ws.Range("E14").Value = '  -2.38%  '
# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.415'
$ws.Range("E15").Value = '  -2.92%  '
# Row 16
$ws.Range("D16").Value = '1.761.90'
$ws.Range("E16").Value = '  -2.15%  '
# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '90.24'
$ws.Range("E17").Value = '  +11.85%  '
# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.00001070'
$ws.Range("E18").Value = '  -1.92%  '
# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06244'
$ws.Range("E19").Value = '  -8.48%  '
# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '1.002'
$ws.Range("E20").Value = '  +0.18%  '
# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '17.35'
$ws.Range("E21").Value = '  -0.91%  '
# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.166'
$ws.Range("E22").Value = '  -2.25%  '
# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.5314'
$ws.Range("E23").Value = '  -2.81%  '
# Row 24
$ws.Range("D24").Value = '27.521.73'
$ws.Range("E24").Value = '  -2.98%  '
# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '11.52'
$ws.Range("E25").Value = '  -2.40%  '
# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.306'
$ws.Range("E26").Value = '  -4.35%  '
# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '20.49'
$ws.Range("E27").Value = '  +0.07%  '
# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '153.09'
$ws.Range("E28").Value = '  -0.60%  '
# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.296'
$ws.Range("E29").Value = '  -2.43%  '
# Row 30
$ws.Range("D30").Value = '1.959.90'
# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '127.41'
$ws.Range("E31").Value = '  -3.61%  '
# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.174'
$ws.Range("E32").Value = '  -6.62%  '
# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.704'
$ws.Range("E33").Value = '  -1.75%  '
# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.09188'
$ws.Range("E34").Value = '  -1.50%  '
# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '3.631'
$ws.Range("E35").Value = '  -9.37%  '
# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '12.56'
$ws.Range("E36").Value = '  +3.79%  '
# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.02307'
$ws.Range("E37").Value = '  -0.68%  '
# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.2149'
# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.06092'
$ws.Range("E39").Value = '  -4.41%  '
# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.6423'
$ws.Range("E40").Value = '  -2.51%  '
# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '5.037'
$ws.Range("E41").Value = '  -2.25%  '
# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.174'
$ws.Range("E42").Value = '  -2.90%  '
# Row 43
$ws.Range("B43").Value = 'FraxShare'
$ws.Range("C43").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '7.940'
$ws.Range("E43").Value = '  -2.36%  '
# Row 44
$ws.Range("B44").Value = 'Frax'
$ws.Range("C44").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.001'
$ws.Range("E44").Value = '  +0.17%  '
# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.385'
$ws.Range("E45").Value = '  -5.00%  '
# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '13.74'
$ws.Range("E46").Value = '  -0.66%  '
# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.5955'
$ws.Range("E47").Value = '  -2.08%  '
# Row 48
$ws.Range("E48").Value = '  -2.22%  '
# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '126.24'
$ws.Range("E49").Value = '  -1.60%  '
# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.977'
$ws.Range("E50").Value = '  -2.73%  '
# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.06877'
$ws.Range("E51").Value = '  -3.00%  '
